$d = $word.ActiveDocument

# The edit rewrites two placeholder paragraphs in the "Draft Class" section:
#   1) The heading "The Draft Class" -> The <curly-quoted "Draft"> Class, split
#      across 5 runs (matching the quoting style used by the sibling headings
#      like The "Ranking"/"League"/"CustomLeague" Class).
#   2) The placeholder body "Add info about the draft class" -> several
#      sentences of real content, split across multiple runs, spilling into a
#      second new paragraph that also carries a lastRenderedPageBreak marker.
#
# Because Range.Text only ever produces a single merged run, the precise
# multi-run / multi-paragraph structure is applied with Range.InsertXML,
# which replaces the full contents of the target paragraph range with the
# supplied OOXML (one or more <w:p> elements).

$headingPara = $null
$bodyPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text
    if ($t -eq "The Draft Class`r") {
        $headingPara = $para
    }
    if ($t -eq "Add info about the draft class`r") {
        $bodyPara = $para
    }
}

if ($headingPara -eq $null) {
    throw "Could not find the `"The Draft Class`" heading paragraph"
}
if ($bodyPara -eq $null) {
    throw "Could not find the `"Add info about the draft class`" body paragraph"
}

$headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">The </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr><w:t>“</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr><w:t>Draft</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr><w:t>”</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Class</w:t></w:r>' +
    '</w:p>'

$bodyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The overall purpose of the draft class is to have the user or computer select from a list of players to fill out their fantasy football team. </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">This class consists of private members that initialize important variables to be used later. These include the positional and round limits, player average salaries, </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">containers for the players, and many more. </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The draft constructor takes in the league members, position limits, round limit, and draft type. </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">In its body, the constructor fills the team player map with the league members, and then chooses which CSV file to open based on the user input. It then resizes the positional limit containers depending on the number of league members. </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Next is the displayTopPlayers function. This displays function overloading because there are two instances. The first instance displays the top ten players available to the user. The second instance allows the user to input the number of players they would like to see, and then outputs that number of players. </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">After this the promptForMorePlayers function allows the user to see a different number of top players if they were not initially satisfied with their choice. This is useful when searching for a player in a certain position. The function makes sure to check that the user inputs are valid for the situation and do not cause any unwanted output. The updatePositionCount function then adjusts the </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">position limit container so no drafter exceeds the limit on any position. </w:t></w:r>' +
    '</w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Next up is the getUserPick function, which allows the user to pick the player they would like to draft. This function receives an input from the user, checks that it is a positive number, and </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">checks that this selection does not exceed the positional limits. If the input passes all of these cases, the player is added to the team’s roster. </w:t></w:r>' +
    '</w:p>'

$headingPara.Range.InsertXML($headingXml)
$bodyPara.Range.InsertXML($bodyXml)

Write-Output "Rewrote the Draft Class heading and body paragraphs"
